$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $targetText) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd("`r") -eq $targetText) {
            return $p
        }
    }
    return $null
}

# 1. Title heading replacement. The same sentence also appears, verbatim, as
#    the bold blurb near the end of the document; the diff changes both
#    occurrences identically, so a document-wide replace handles both.
$d.Content.Find.Execute(
    "Play Joker Lanterns Free - A Spooky and Festive Slot Game", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Joker Lanterns Free - Classic Halloween Slot Game", 2) | Out-Null

# 2. "What we like" bullet list
#    a) reword the first bullet
$d.Content.Find.Execute(
    "Halloween theme with classic and modern elements", $true, $false, $false, $false, $false,
    $true, 1, $false, "Classic Halloween theme with a festive and spooky feel", 2) | Out-Null

#    b) drop the closing "Special symbols..." bullet entirely
$deadPara = Get-ParagraphByText $d "Special symbols that contribute to the excitement of the game"
if ($deadPara -ne $null) {
    $deadPara.Range.Delete()
}

#    c) insert a brand-new bullet right after the reworded first one
$anchor = Get-ParagraphByText $d "Classic Halloween theme with a festive and spooky feel"
if ($anchor -ne $null) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()
    $newPara.Range.Text = "Variety of symbols that blend classic and modern elements"
}

#    d) reword the (now-last) "Free spins mode" bullet
$d.Content.Find.Execute(
    "Free spins mode with rich multipliers and bonus spins", $true, $false, $false, $false, $false,
    $true, 1, $false, "Free spins mode with rich multipliers and additional spins", 2) | Out-Null

# 3. "What we don't like" bullet list
$d.Content.Find.Execute(
    "May be too spooky for some players", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited availability of special symbols and jackpots", 2) | Out-Null

$d.Content.Find.Execute(
    "Limited variety of symbols compared to other modern slot games", $true, $false, $false, $false, $false,
    $true, 1, $false, "May not appeal to players looking for a completely modern slot game", 2) | Out-Null

# 4. Italic blurb at the very end of the document
$d.Content.Find.Execute(
    "Read our review of Joker Lanterns, a classic Halloween-themed slot game with a modern twist. Play for free and win big prizes with special symbols.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Joker Lanterns, a classic-themed slot game with a modern twist. Play for free!", 2) | Out-Null

Write-Output "edits applied"
Write-Output ("Paragraph count: " + $d.Paragraphs.Count)

